$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row above the current row 23 ("Sin especificar" /
# 2023-10-11 record). Everything from the old row 23 downward shifts down
# by one row (old row 23 -> new row 24, ..., old row 120 -> new row 121).
$ws.Rows(23).Insert()

# Populate the newly inserted row 23 with the new weekly record.
$ws.Cells.Item(23, 1).Value = 5
$ws.Cells.Item(23, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(23, 3).Value = "Maule"
$ws.Cells.Item(23, 4).Value = 45222
$ws.Cells.Item(23, 5).Value = 7
$ws.Cells.Item(23, 6).Value = 300000000
$ws.Cells.Item(23, 7).Value = "Espárragos"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 4000
$ws.Cells.Item(23, 11).Value = 1100
$ws.Cells.Item(23, 12).Value = 1200
$ws.Cells.Item(23, 13).Value = 1150
$ws.Cells.Item(23, 14).Value = "$/kilo"
$ws.Cells.Item(23, 15).Value = "Provincia de Linares"
$ws.Cells.Item(23, 16).Value = 1150
$ws.Cells.Item(23, 17).Value = 1
$ws.Cells.Item(23, 18).Value = "Hortaliza"
